$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$colB = @(1.47377192904213, 1.337329118349317, 1.254203563237297, 1.220490691248528, 1.214902383359401, 1.253748248768829, 1.426589975161392, 1.770815637081853, 2.027132682713386, 2.144523614702337, 2.18909336101575, 2.179489251970324, 2.148188047637461, 2.129030386430202, 2.019477047381372, 1.952474051356717, 1.91400999548506, 1.900999409206065, 1.959598926762567, 2.157378791399083, 2.287318785944876, 2.217904322084337, 1.956377593683044, 1.677108060023556)
for ($i = 0; $i -lt $colB.Length; $i++) {
    $ws.Cells.Item(2 + $i, 2).Value = $colB[$i]
}

$colC = @(0.5293168308653549, 0.475428262252592, 0.442588014119849, 0.4292661127416295, 0.4270576484061053, 0.4424081061930565, 0.510683759037363, 0.6466106443111812, 0.7478326059448932, 0.7942011876249921, 0.8118079951660206, 0.8080138905031617, 0.7956487377960002, 0.7880810308923856, 0.7448089382040166, 0.718346634311672, 0.7031564408072768, 0.6980184419923603, 0.7211604514765213, 0.7992793656060257, 0.850615370589253, 0.8231901174949598, 0.7198882521738597, 0.6096089300387462)
for ($i = 0; $i -lt $colC.Length; $i++) {
    $ws.Cells.Item(2 + $i, 3).Value = $colC[$i]
}

$colE = @(0.3836592374743972, 0.3747380058281706, 0.3694211161247694, 0.367294877288316, 0.3669442595296957, 0.3693922771833016, 0.380549825759104, 0.4037069028737861, 0.4215034712468366, 0.4297708156572924, 0.4329261747138133, 0.4322455129933473, 0.4300299138868695, 0.4286760121711666, 0.4209666353843247, 0.4162811449002461, 0.4136023246621292, 0.4126980954053181, 0.4167782517734366, 0.4306800186496815, 0.4399095819350691, 0.4349704105765611, 0.4165534632518231, 0.3973051252951763)
for ($i = 0; $i -lt $colE.Length; $i++) {
    $ws.Cells.Item(2 + $i, 5).Value = $colE[$i]
}

$colF = @(2.891181040996173, 2.830483127117958, 2.795192298004082, 2.781305224382265, 2.779029045378607, 2.795003015135634, 2.869839950710173, 3.032443416710578, 3.16180208601665, 3.222851098389384, 3.24628955965855, 3.241227352226645, 3.224772949347596, 3.214736005907696, 3.157857038181334, 3.123530528601577, 3.103993968338074, 3.097414688906611, 3.127163179189068, 3.22959727737782, 3.298413923597224, 3.261512828257679, 3.125520240767287, 2.986734398785359)
for ($i = 0; $i -lt $colF.Length; $i++) {
    $ws.Cells.Item(2 + $i, 6).Value = $colF[$i]
}

$colG = @(0.002454686167697295, 0.002460261272735173, 0.002463860007106984, 0.002465370836960481, 0.002465624390235558, 0.002463880203044178, 0.002456572115793153, 0.002443626900686007, 0.002434950600382235, 0.002431182525044089, 0.002429781196578301, 0.002430081863420699, 0.002431066725492322, 0.00243167330605794, 0.002435200438040764, 0.002437409909131132, 0.002438697578806326, 0.002439136458540214, 0.002437172965483912, 0.002430776755269512, 0.002426745373692809, 0.002428883423372108, 0.002437280033442629, 0.002446981622975619)
for ($i = 0; $i -lt $colG.Length; $i++) {
    $ws.Cells.Item(2 + $i, 7).Value = $colG[$i]
}

$colI = @(0.8138053962467424, 0.8091280925357935, 0.8069250478298287, 0.8061942285513695, 0.8060829230661568, 0.806914517472805, 0.8120530878604768, 0.8274949545811126, 0.8421943742457358, 0.849628299037434, 0.8525522427134575, 0.8519176534950361, 0.8498666632107756, 0.848624594829289, 0.8417237063553671, 0.8376826051880641, 0.8354285021427188, 0.8346773200944781, 0.8381055081932089, 0.8504661222284966, 0.8591798064633451, 0.8544705349502522, 0.837914098374128, 0.8227345738192327)
for ($i = 0; $i -lt $colI.Length; $i++) {
    $ws.Cells.Item(2 + $i, 9).Value = $colI[$i]
}

$colJ = @(0.09877265073556174, 0.09688935948992849, 0.09581831636496929, 0.09540309984286921, 0.09533543037960612, 0.09581263087387981, 0.09810545992098696, 0.1032883742768291, 0.1075290961652584, 0.1095554400434438, 0.1103369935860954, 0.1101680356817099, 0.1096194525171441, 0.109285288715661, 0.1073986466478942, 0.1062663233032453, 0.1056241768684174, 0.1054083189435531, 0.1063859137784817, 0.1097801967940555, 0.1120815628693919, 0.1108456048073805, 0.106331819404204, 0.1018112867991476)
for ($i = 0; $i -lt $colJ.Length; $i++) {
    $ws.Cells.Item(2 + $i, 10).Value = $colJ[$i]
}

$colN = @(1.151298788324876, 1.170788999005773, 1.183336827568893, 1.188595778362267, 1.189477805157038, 1.183407162637646, 1.1578982190588, 1.112501839705843, 1.081994084874653, 1.068738472568658, 1.063808940242783, 1.064866587782413, 1.068331108366529, 1.070464975293525, 1.082872944826009, 1.090644670899025, 1.095173322179377, 1.096716686452956, 1.089811292616561, 1.067311045386825, 1.053131160249031, 1.06065097519282, 1.0901878742799, 1.124284968391359)
for ($i = 0; $i -lt $colN.Length; $i++) {
    $ws.Cells.Item(2 + $i, 14).Value = $colN[$i]
}
